$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.367.55"
$ws.Range("E2").Value = "  -4.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.563.87"
$ws.Range("E3").Value = "  -4.21%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.29"
$ws.Range("E6").Value = "  -2.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3693"
$ws.Range("E7").Value = "  -2.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.33"
$ws.Range("E8").Value = "  -1.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3390"
$ws.Range("E9").Value = "  -2.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.165"
$ws.Range("E10").Value = "  -3.61%  "
$ws.Range("E11").Value = "  -4.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("E13").Value = "  -2.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.056"
$ws.Range("E14").Value = "  -3.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.911"
$ws.Range("E15").Value = "  -4.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.562.79"
$ws.Range("E16").Value = "  -3.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001126"
$ws.Range("E17").Value = "  -5.88%  "
$ws.Range("E18").Value = "  -4.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06732"
$ws.Range("E19").Value = "  -3.13%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.246"
$ws.Range("E21").Value = "  -5.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.51"
$ws.Range("E22").Value = "  -4.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5301"
$ws.Range("E23").Value = "  -6.68%  "
$ws.Range("E24").Value = "  -3.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "22.362.34"
$ws.Range("E25").Value = "  -4.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.368"
$ws.Range("E26").Value = "  -2.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.818"
$ws.Range("E27").Value = "  -4.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.10"
$ws.Range("E28").Value = "  -3.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "146.69"
$ws.Range("E29").Value = "  -2.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.974"
$ws.Range("E30").Value = "  -3.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.73"
$ws.Range("E31").Value = "  -4.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.738.30"
$ws.Range("E32").Value = "  -3.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.017"
$ws.Range("E33").Value = "  +3.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.193"
$ws.Range("E34").Value = "  -8.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.015"
$ws.Range("E35").Value = "  -5.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.05"
$ws.Range("E36").Value = "  -9.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08459"
$ws.Range("E37").Value = "  -3.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02538"
$ws.Range("E38").Value = "  -5.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2320"
$ws.Range("E39").Value = "  -3.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.515"
$ws.Range("E40").Value = "  -5.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06455"
$ws.Range("E41").Value = "  -5.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.276"
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("E43").Value = "  -8.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6338"
$ws.Range("E44").Value = "  -6.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.20"
$ws.Range("E45").Value = "  -7.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("E47").Value = "  -5.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.753"
$ws.Range("E48").Value = "  -3.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.102"
$ws.Range("E49").Value = "  -5.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.264"
$ws.Range("E50").Value = "  +4.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "124.78"
$ws.Range("E51").Value = "  -1.49%  "
